$d = $word.ActiveDocument

# The site footer block ("Ver no Jupiter Salvar em pdf Salvar em docx" +
# the "(c) 2020 ... Jekyll ..." line, plus the blank paragraph that
# separates them from the requirements text above) is being dropped from
# this rebuilt page. Find those paragraphs by their text and remove them
# (each Range.Delete also removes that paragraph's own paragraph mark).

$jupiterText = "Ver no Jupiter Salvar em pdf Salvar em docx"
$copyrightSnippet = "Powered by Jekyll and Github pages"

$jupiterIndex = -1
$copyrightIndex = -1

$paras = $d.Paragraphs
$count = $paras.Count

for ($i = 1; $i -le $count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -like "*$jupiterText*") {
        $jupiterIndex = $i
    }
    if ($t -like "*$copyrightSnippet*") {
        $copyrightIndex = $i
    }
}

if ($jupiterIndex -gt 0 -and $copyrightIndex -eq ($jupiterIndex + 1)) {
    $blankIndex = $jupiterIndex - 1
    $blankText = $d.Paragraphs.Item($blankIndex).Range.Text.Trim()

    # Delete from the highest index down so lower indices stay valid.
    $d.Paragraphs.Item($copyrightIndex).Range.Delete()
    $d.Paragraphs.Item($jupiterIndex).Range.Delete()
    if ($blankText -eq "") {
        $d.Paragraphs.Item($blankIndex).Range.Delete()
    }
}
